$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Default footer (footer1.xml): Pearson logo id=1, image1.png -> image2.png ---
$ftr1 = $sec.Footers.Item(1)
$pearson1 = $ftr1.Range.InlineShapes.Item(1)
$shape1 = $pearson1.ConvertToShape()
$shape1.Name = "image2.png"
$shape1.ConvertToInlineShape() | Out-Null

# --- First-page footer (footer2.xml): Pearson logo id=2, image1.png -> image2.png ---
$ftr2 = $sec.Footers.Item(2)
$pearson2 = $ftr2.Range.InlineShapes.Item(1)
$shape2 = $pearson2.ConvertToShape()
$shape2.Name = "image2.png"
$shape2.ConvertToInlineShape() | Out-Null

# --- First-page header (header2.xml): BTEC logo id=3, image2.jpg -> image1.jpg ---
$hdr2 = $sec.Headers.Item(2)
$btec = $hdr2.Range.InlineShapes.Item(1)
$shape3 = $btec.ConvertToShape()
$shape3.Name = "image1.jpg"
$shape3.ConvertToInlineShape() | Out-Null

Write-Host "Renamed inline picture shapes in footer1, footer2 and header2."
